$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(35.04526029724706,34.53409286247759,34.22957440854315,34.1079653765856,34.08792621246343,34.22792411090392,34.86715450508282,36.18857139052859,37.19224761578318,37.65413556258149,37.82965173877867,37.79182651893419,37.66856404745985,37.59313703441692,37.16215574826833,36.89901211502369,36.74816781877733,36.69718669805045,36.92697265643984,37.70475390389535,38.2165608442696,37.94313139914119,36.91433031190648,35.82470126194112)
$CValues = @(27.35262150382964,26.87710213571744,26.59171927233644,26.47722315146734,26.45832396236959,26.59016766430281,27.187381992481,28.40437248930845,29.31781586436463,29.73575847935745,29.89422380956509,29.86008888023382,29.74879234352633,29.68064170312367,29.29053816070162,29.05173371772509,28.91461533645909,28.86823445186216,29.07713155726454,29.78147856693077,30.24290874406367,29.99658149598653,29.06564864243126,28.07112958950661)
$DValues = @(15.53448649421699,15.50755928449795,15.49530139205939,15.49137795353404,15.49079108509568,15.49524414478268,15.52431284209825,15.61535464813553,15.70309772190123,15.74754666666588,15.76502903063703,15.76123499197296,15.748971925633,15.74154510034351,15.70028400999751,15.67613219140451,15.66266762863669,15.65818211507399,15.67865901676003,15.75255624693321,15.80464275439151,15.77649709841456,15.67751532871493,15.58706706259622)
$EValues = @(17.01119759690323,16.98736698065072,16.97741894415568,16.97453844434369,16.97413086311037,16.97737535214791,17.00200621405633,17.08763480313596,17.17347077674012,17.21751505252582,17.23491179566497,17.23113317315969,17.21893190630384,17.21155177170867,17.17069291634179,17.14690717294,17.13369667417371,17.12930462198536,17.14939051837343,17.22249623219596,17.27446003610014,17.2463434233022,17.14826635186955,17.06045273203367)
$GValues = @(3.8249302817056,3.832624984625362,3.837571473140339,3.839643368931549,3.839990808177397,3.837599187578047,3.827537576315756,3.809550980913862,3.797376163658723,3.792058227211552,3.790075753623887,3.790501328320558,3.791894502583251,3.792751928419375,3.797728106367651,3.800837019093418,3.80264595661101,3.803262011277987,3.800503922869752,3.791484446798701,3.785772033646161,3.788804301501576,3.800654448419334,3.814232504093363)
$JValues = @(9.765152013474401,9.789423655969607,9.806154120453375,9.813430034923137,9.814665821741928,9.806250392986177,9.773140691349949,9.722779462338936,9.694751337195395,9.683968838359784,9.680170369256119,9.680975751654417,9.683650621945654,9.685326174523547,9.695495736613857,9.702239588996957,9.706303576517881,9.707711313075237,9.701502524865942,9.682857209469308,9.672331243258265,9.677796716196545,9.70183516966919,9.734834553130227)
$NValues = @(24.77426605190089,24.67782422435335,24.61962592075575,24.59616704376069,24.59228730137962,24.61930850207267,24.74080277601543,24.98705142078509,25.1727435409501,25.258229181801,25.29074302750703,25.28373431782254,25.2609012842929,25.24693381632131,25.16717695674197,25.11850588713025,25.0906071789694,25.08117767649545,25.12367713506979,25.26760406649637,25.36249625647379,25.31177597314553,25.12133895528479,24.91958776192925)

$startRow = 2

for ($i = 0; $i -lt $BValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value2 = $BValues[$i]
}

for ($i = 0; $i -lt $CValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value2 = $CValues[$i]
}

for ($i = 0; $i -lt $DValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value2 = $DValues[$i]
}

for ($i = 0; $i -lt $EValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value2 = $EValues[$i]
}

for ($i = 0; $i -lt $GValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value2 = $GValues[$i]
}

for ($i = 0; $i -lt $JValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 10).Value2 = $JValues[$i]
}

for ($i = 0; $i -lt $NValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 14).Value2 = $NValues[$i]
}
